# "Changing layout to expand on RQs"
#
# Slide 4 (sldId 263, "Research Questions") is trimmed down to just the
# first research-question textbox (reworded) plus the surrounding
# decorative shapes that stay; every other shape that used to hold
# Research Question 1/2/3 boxes, the "missing nodes" question and the
# "blockchain" question is removed.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Reword the remaining "Research Question" textbox (shape id 8) ----
$sh8 = Get-ShapeById $s 8
$tf8 = $sh8.TextFrame
$tr8 = $tf8.TextRange

$run2 = $tr8.Runs(2, 1)
$run2.Text = "can FL algorithms or tools be designed to operate specifically on Edge devices to minimize computational and costs"

$run3 = $tr8.Runs(3, 1)
$run3.Text = ". The models should be adaptable to a range of data sources and types."

# --- Delete the shapes that were removed from the slide ---------------
# id 3  - TextBox 2 ("When thinking about robust and resilient FL systems, can blockchain ")
# id 4  - Rectangle: Rounded Corners 3 (frame around shape 3)
# id 10 - TextBox 9 ("In any large distributed network ... Can a FL architecture ...")
# id 14 - TextBox 13 ("Research Question 1")
# id 15 - Rectangle: Rounded Corners 14 (frame around shape 10)
# id 17 - Rectangle: Rounded Corners 16 (frame around shape 19)
# id 18 - TextBox 17 ("Research Question 3")
# id 19 - TextBox 18 ("Can we design efficient aggregation algorithms for FL ...")
$idsToDelete = @(3, 4, 10, 14, 15, 17, 18, 19)
foreach ($id in $idsToDelete) {
    $sh = Get-ShapeById $s $id
    if ($sh -ne $null) {
        $sh.Delete()
    }
}
